$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need an explicit Text number format
# so Excel does not coerce them (and strip things like trailing zeros).
$textCells = @("D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "22.392.48"
$ws.Range("E2").Value = "  -4.73%  "
$ws.Range("D3").Value = "1.571.55"
$ws.Range("E3").Value = "  -4.70%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "290.89"
$ws.Range("E6").Value = "  -3.06%  "
$ws.Range("D7").Value = "0.3678"
$ws.Range("E7").Value = "  -3.09%  "
$ws.Range("D8").Value = "49.56"
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").Value = "0.3372"
$ws.Range("E9").Value = "  -5.44%  "
$ws.Range("D10").Value = "1.166"
$ws.Range("E10").Value = "  -4.77%  "
$ws.Range("D11").Value = "0.07563"
$ws.Range("E11").Value = "  -6.72%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "21.10"
$ws.Range("E13").Value = "  -4.47%  "
$ws.Range("D14").Value = "6.054"
$ws.Range("E14").Value = "  -5.65%  "
$ws.Range("D15").Value = "6.853"
$ws.Range("E15").Value = "  -7.57%  "
$ws.Range("D16").Value = "0.00001137"
$ws.Range("E16").Value = "  -5.44%  "
$ws.Range("D17").Value = "1.573.30"
$ws.Range("E17").Value = "  -5.01%  "
$ws.Range("D18").Value = "89.30"
$ws.Range("E18").Value = "  -8.21%  "
$ws.Range("D19").Value = "0.06726"
$ws.Range("E19").Value = "  -3.70%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "6.227"
$ws.Range("E21").Value = "  -8.07%  "
$ws.Range("D22").Value = "16.31"
$ws.Range("E22").Value = "  -6.59%  "
$ws.Range("D23").Value = "11.93"
$ws.Range("E23").Value = "  -4.96%  "
$ws.Range("D24").Value = "22.415.01"
$ws.Range("E24").Value = "  -4.73%  "
$ws.Range("D25").Value = "2.417"
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("D26").Value = "2.958"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").Value = "19.80"
$ws.Range("E27").Value = "  -5.67%  "
$ws.Range("D28").Value = "145.92"
$ws.Range("E28").Value = "  -4.54%  "
$ws.Range("E29").Value = "  -5.78%  "
$ws.Range("E30").Value = "  -6.07%  "
$ws.Range("D31").Value = "1.749.97"
$ws.Range("E31").Value = "  -4.78%  "
$ws.Range("D32").Value = "6.243"
$ws.Range("E32").Value = "  -10.11%  "
$ws.Range("D33").Value = "1.974"
$ws.Range("E33").Value = "  -7.88%  "
$ws.Range("D34").Value = "0.9812"
$ws.Range("E34").Value = "  -4.77%  "
$ws.Range("D35").Value = "10.38"
$ws.Range("E35").Value = "  -12.82%  "
$ws.Range("D36").Value = "0.08445"
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("D37").Value = "0.02541"
$ws.Range("E37").Value = "  -6.95%  "
$ws.Range("D38").Value = "0.2296"
$ws.Range("E38").Value = "  -6.53%  "
$ws.Range("D39").Value = "0.06504"
$ws.Range("E39").Value = "  -5.13%  "
$ws.Range("D40").Value = "5.489"
$ws.Range("E40").Value = "  -8.13%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "11.75"
$ws.Range("E41").Value = "  -12.66%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.257"
$ws.Range("E42").Value = "  -4.93%  "
$ws.Range("D43").Value = "0.6384"
$ws.Range("E43").Value = "  -7.90%  "
$ws.Range("D44").Value = "14.49"
$ws.Range("E44").Value = "  -7.57%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "0.5994"
$ws.Range("E46").Value = "  -7.07%  "
$ws.Range("D47").Value = "3.773"
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("D48").Value = "2.113"
$ws.Range("D49").Value = "120.88"
$ws.Range("E49").Value = "  -5.73%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.07281"
$ws.Range("E50").Value = "  -6.69%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "1.191"
$ws.Range("E51").Value = "  +0.31%  "
